$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Set the new value entered for Cliente (row 7) in column N (iteration 6 payment)
$ws.Range("N7").Value = 1

# Recalculate dependent shared formulas across the sheet
$excel.Calculate()

# Reflect the active cell selection on the frozen-pane sheet view (bottom-right pane -> N7)
$ws.Range("N7").Select()
